# Update the Japanese "Planned Corporate Expansion" challenge bullets in
# the Adatum marketing report to the revised wording.

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false, $true,
        $wdFindContinue, $false, $replace, $wdReplaceAll)
    Write-Output $ok
}

# "限定的なブランド認識と認識" (bullet heading)
Replace-Text `
    "限定的なブランド認識と認識" `
    "ブランドの認知度と知名度の不足"

# Bullet 1 body text
Replace-Text `
    ": これらの新しい市場で可視性を達成することは主要なハードルであり、Adatumのブランドプレゼンスをゼロから構築するための堅牢なマーケティング努力が必要です。" `
    ": 新しい市場における認知度の向上は重要な課題であり、Adatum のブランド プレゼンスをゼロから構築するにはマーケティング活動を強化する必要があります。"

# Bullet 2 ("激しい競争") body text
Replace-Text `
    ": カナダのクラウド サービス部門は、多くのプレーヤーと激しく競争しています。" `
    ": カナダのクラウド サービス業界は競争が激しく、多くの競合他社が存在します。"

# Bullet 3 ("多様な顧客の好みと期待") body text
Replace-Text `
    ": これらの市場のさまざまな要求に合わせて製品とマーケティングを調整することは、地元の企業や消費者との共感を得るために不可欠です。" `
    ": カナダ市場の多様な需要に合わせて製品やマーケティングをカスタマイズし、地元の企業や消費者の共感を得る必要があります。"

# Bullet 4 ("規制とコンプライアンスの課題") body text
Replace-Text `
    ": Adatum は、地域の個別のデータ プライバシー、セキュリティ、運用に関する規制をナビゲートするという複雑なタスクに直面し、勤勉なコンプライアンスの取り組みを必要とします。" `
    ": Adatum は、地域特有のデータ プライバシー、セキュリティ、運用規制を遵守するという複雑な課題に直面しており、慎重なコンプライアンス対応が求められています。"

# Bullet 5 ("運用と物流の複雑さ") body text
Replace-Text `
    ": 効率的なリージョン間の運用を確立することは、特に高いサービス レベルの維持と地理的な場所にまたがるデータ センターの管理において、物流上の課題を提示します。" `
    ": 地域を越えて効率的な運用を確立するには、高いサービス レベルの維持や地理的に離れているデータ センターの管理などの物流に関する課題を解消する必要があります。"
